# Add files via upload
# Re-upload of the broadcast schedule: adds the station logo hyperlink/image
# cell (column G) to rows that were missing it, and restores the normal
# (unscrolled) sheet view/selection that Excel saves on close.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$logoUrl = "https://static.wixstatic.com/media/c39cce_e474bc38cd9d46a5b7af4e985ee1892c~mv2.png/v1/fill/w_150,h_153,al_c,q_85,usm_0.66_1.00_0.01,enc_auto/logo%20radio%20sol%20-%20color%201.png"

# Cells in column G that need the logo hyperlink added, in authoring order.
$targetCells = @("G19","G23","G24","G29","G31","G33","G34","G35","G38","G40","G41","G45","G46","G47","G67","G66","G69","G71","G73","G74","G76","G77","G79","G80","G60")

# The existing hyperlink style already used throughout column G (e.g. G12).
$hyperlinkStyle = $ws.Range("G12").Style

foreach ($cellRef in $targetCells) {
    $cell = $ws.Range($cellRef)
    $ws.Hyperlinks.Add($cell, $logoUrl)
    # Hyperlinks.Add() stamps its own (duplicate) hyperlink style onto the
    # cell; reapply the sheet's existing hyperlink style so the cell matches
    # the style already shared by every other logo cell in column G.
    $cell.Style = $hyperlinkStyle
}

# Restore the plain/unscrolled view with G60 selected (matches the state the
# workbook was saved in), clearing the old frozen "topLeftCell"/selection.
$ws.Range("G60").Select()

Write-Output "done"
